$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27 - PeakValue: add formula to Details (C27), set Unit (D27) to '% time^-1^'
$c27Text = @'
It is the accumulated number of seeds germinated at the point on the germination curve at which the rate of germination starts to decrease. It is computed as the maximum quotient obtained by dividing successive cumulative germination values by the relevant incubation time.
$$PV = \max\left ( \frac{G_{1}}{T_{1}},\frac{G_{2}}{T_{2}},\cdots \frac{G_{k}}{T_{k}} \right )$$
Where, Where, $T_{i}$ is the time from the start of the experiment to the
$i$th interval, $G_{i}$ is the cumulative germination percentage in the $i$th time interval and $k$ is the total number of time intervals.
'@
$ws.Range("C27").Value = $c27Text
$ws.Range("D27").Value = '\% time^-1^'
$ws.Range("D27").Style = $ws.Range("D12").Style

# Row 28 - GermValue (Czabator): extend Details (C28) with GV_mod sentence; Reference (F28) gets extra citation
$c28Text = @'
It is computed as follows:
$$GV = PV \times MDG$$
Where, $PV$ is the peak value and $MDG$ is the mean daily germination percentage from the onset of germination.
It can also be computed for other time intervals of successive germination counts, by replacing $MDG$ with the mean germination percentage per unit time ($\overline{GP}$).
$GV$ value can be modified ($GV_{mod}$), to consider the entire duration from the beginning of the test instead of just from the onset of germination.
'@
$ws.Range("C28").Value = $c28Text
$ws.Range("F28").Value = '[@czabator_germination_1962; @brown_representing_1988]'

# New yellow highlight style (fillId=2) applied to D28, D29, D20 (Unit column cells)
$ws.Range("D28").Interior.Color = 65535
$ws.Range("D28").HorizontalAlignment = -4131
$ws.Range("D28").VerticalAlignment = -4160
$ws.Range("D28").WrapText = $false

$ws.Range("D29").Interior.Color = 65535
$ws.Range("D29").HorizontalAlignment = -4131
$ws.Range("D29").VerticalAlignment = -4160
$ws.Range("D29").WrapText = $false

$ws.Range("D20").Interior.Color = 65535
$ws.Range("D20").HorizontalAlignment = -4131
$ws.Range("D20").VerticalAlignment = -4160
$ws.Range("D20").WrapText = $false

# Row heights to fit the new, longer content
$ws.Rows.Item(27).RowHeight = 135
$ws.Rows.Item(28).RowHeight = 120

Write-Output "done"
